$wb = $excel.ActiveWorkbook

# --- Sheet "good" (sheet1): insert a new column B for "en_comments" ---
$ws1 = $wb.Worksheets.Item("good")
$ws1.Columns("B:B").Insert()
$ws1.Range("B1").Value = "en_comments"
$ws1.Columns("B:B").ColumnWidth = 11.92

# --- Sheet "missing_translation_cell" (sheet2): insert two new columns B,C
#     for "en_comments" / "de_comments" ---
$ws2 = $wb.Worksheets.Item("missing_translation_cell")
$ws2.Columns("B:C").Insert()
$ws2.Range("B1").Value = "en_comments"
$ws2.Range("C1").Value = "de_comments"

# New comment values, in the order that reproduces the target shared-string table
$ws2.Range("G4").Value = "Knoten 2"
$ws2.Range("C2").Value = "List Missing English Comment"
$ws2.Range("B4").Value = "Node Missing German Comment"

# Fill out the rest of new row 4 (List 4 / List de / Node 2), matching row2/row3 style
$ws2.Range("D4").Value = "List 4"
$ws2.Range("E4").Value = "List de"
$ws2.Range("F4").Value = "Node 2"
$ws2.Range("D4:G4").Font.Color = 0

$ws2.Columns("B:B").ColumnWidth = 27.59
$ws2.Columns("C:C").ColumnWidth = 21.75

# --- Selection / active sheet bookkeeping ---
$ws2.Range("D5").Select() | Out-Null
$ws1.Range("D12").Select() | Out-Null

Write-Output "done"
